$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "119 methods" -> "117 methods" (shared header label used by A22 and A142)
$ws.Range("A22").Value = "117 methods"
$ws.Range("A142").Value = "117 methods"

# 2) Corrected aggregate row values for the "NbMethodsCallingMe"/other columns
#    (Sum/Average/StdDev/Variance recomputed over the 117 non-blank method rows)

# Sum row (143): only G (NbMethodsCallingMe) sum changes
$ws.Range("G143").Value = 103

# Average row (144)
$ws.Range("B144").Value = 8.632478632478632
$ws.Range("C144").Value = 0.6495726495726496
$ws.Range("D144").Value = 1.3162393162393162
$ws.Range("E144").Value = 1.8376068376068375
$ws.Range("F144").Value = 1.7521367521367521
$ws.Range("G144").Value = 0.88034188034188032

# Standard deviation row (147)
$ws.Range("B147").Value = 8.978224645663806
$ws.Range("C147").Value = 2.369017392742232
$ws.Range("D147").Value = 1.083280650786657
$ws.Range("E147").Value = 1.9740074427126335
$ws.Range("F147").Value = 2.008128707124067
$ws.Range("G147").Value = 0.76418095546178155

# H147 picks up a new number-formatted variant of its existing style (fontId 26 /
# fillId 38 / borderId 10) -> applying a 2-decimal number format on the still-empty
# cell reproduces the new cellXfs entry exactly.
$ws.Range("H147").NumberFormat = "0.00"

# Variance row (148)
$ws.Range("B148").Value = 80.608517788004974
$ws.Range("C148").Value = 5.6122434071152023
$ws.Range("D148").Value = 1.1734969683687633
$ws.Range("E148").Value = 3.8967053838848709
$ws.Range("F148").Value = 4.0325809043757763
$ws.Range("G148").Value = 0.58397253269048144

# 3) Update the view's selection/scroll state to match the saved window position
$excel.ActiveWindow.ScrollRow = 137
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C162").Select() | Out-Null
